$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: mark existing Clean-C rows (column G) with "X" ---
$gRows = @(172,185,199,220,239,277,287,296,298,299,300,301,302,317,318,320,321,325,328,329,330,331,333,335,339,340,344,346,348,349,353,362,364,372,373)
foreach ($r in $gRows) {
    $ws.Cells.Item($r, 7).Value = "X"
}

# --- Step 2: row 199 gets a Clean-Code-Comment (column M) ---
$ws.Cells.Item(199, 13).Value = 'Nur EINE Stufe pro Methode (eine Abstraktionsebene)'

# --- Step 3: append 27 new clean-code rule rows (391-417) ---
# Row 391
$ws.Cells.Item(391, 2).Value = 'Inheritance'
$ws.Cells.Item(391, 4).Value = 'Basisklassen sollten nicht von abgeleiteten Klassen abhängen'
$ws.Cells.Item(391, 7).Value = 'X'

# Row 392
$ws.Cells.Item(392, 2).Value = 'Classes'
$ws.Cells.Item(392, 4).Value = 'Möglichst kleine Interfaces'
$ws.Cells.Item(392, 7).Value = 'X'

# Row 393
$ws.Cells.Item(393, 2).Value = 'Unused'
$ws.Cells.Item(393, 4).Value = 'Kein ungenutzter Code'
$ws.Cells.Item(393, 7).Value = 'X'
$ws.Cells.Item(393, 13).Value = 'NONAUTOMATIC'

# Row 394
$ws.Cells.Item(394, 2).Value = 'Structure'
$ws.Cells.Item(394, 4).Value = 'Geringe vertikale Gültigkeit von Variablen und Methoden'
$ws.Cells.Item(394, 7).Value = 'X'

# Row 395
$ws.Cells.Item(395, 2).Value = 'Naming'
$ws.Cells.Item(395, 4).Value = 'konsistente Benamung'
$ws.Cells.Item(395, 7).Value = 'X'

# Row 396
$ws.Cells.Item(396, 2).Value = 'Coupling'
$ws.Cells.Item(396, 4).Value = 'Keine überflüssige Kopplung'
$ws.Cells.Item(396, 7).Value = 'X'
$ws.Cells.Item(396, 13).Value = 'NONAUTOMATIC'

# Row 397
$ws.Cells.Item(397, 2).Value = 'Coupling'
$ws.Cells.Item(397, 4).Value = 'Möglichst kein "Funktionsneid"'
$ws.Cells.Item(397, 7).Value = 'X'
$ws.Cells.Item(397, 13).Value = 'NONAUTOMATIC'

# Row 398
$ws.Cells.Item(398, 2).Value = 'Naming'
$ws.Cells.Item(398, 4).Value = 'sprechende Namen für alles - keine Ungarische Notation usw.'
$ws.Cells.Item(398, 7).Value = 'X'
$ws.Cells.Item(398, 13).Value = 'NONAUTOMATIC'

# Row 399
$ws.Cells.Item(399, 2).Value = 'Naming'
$ws.Cells.Item(399, 3).Value = 'Structure'
$ws.Cells.Item(399, 4).Value = 'Zuständigkeiten sinnvoll zuordnen und durch entsprechende Namen kennzeichnen'
$ws.Cells.Item(399, 7).Value = 'X'
$ws.Cells.Item(399, 13).Value = 'NONAUTOMATIC'

# Row 400
$ws.Cells.Item(400, 2).Value = 'Static'
$ws.Cells.Item(400, 3).Value = 'Methods'
$ws.Cells.Item(400, 4).Value = 'Nur Methoden die nie polymorphes Verhalten haben können als Static definieren'
$ws.Cells.Item(400, 7).Value = 'X'
$ws.Cells.Item(400, 13).Value = 'NONAUTOMATIC'

# Row 401
$ws.Cells.Item(401, 2).Value = 'Naming'
$ws.Cells.Item(401, 3).Value = 'Variables'
$ws.Cells.Item(401, 4).Value = 'Zwischenergebnisse in Aussagekräftig benannten Variablen ablegen'
$ws.Cells.Item(401, 7).Value = 'X'
$ws.Cells.Item(401, 13).Value = 'NONAUTOMATIC'

# Row 402
$ws.Cells.Item(402, 2).Value = 'Coding'
$ws.Cells.Item(402, 4).Value = 'Funktionsweise des implementierten Algorithmus verstehen'
$ws.Cells.Item(402, 7).Value = 'X'
$ws.Cells.Item(402, 13).Value = 'NONAUTOMATIC'

# Row 403
$ws.Cells.Item(403, 2).Value = 'Coupling'
$ws.Cells.Item(403, 4).Value = 'Logische Abhängigkeiten in Physischen Abhängigkeiten festschreiben'
$ws.Cells.Item(403, 7).Value = 'X'
$ws.Cells.Item(403, 13).Value = 'NONAUTOMATIC'

# Row 404
$ws.Cells.Item(404, 2).Value = 'Conditionals'
$ws.Cells.Item(404, 4).Value = 'Switch-Anweisungen vermeiden - stattdessen polymorphismus'
$ws.Cells.Item(404, 7).Value = 'X'

# Row 405
$ws.Cells.Item(405, 2).Value = 'Coding'
$ws.Cells.Item(405, 4).Value = 'Präzise implementieren'
$ws.Cells.Item(405, 7).Value = 'X'
$ws.Cells.Item(405, 13).Value = 'NONAUTOMATIC'

# Row 406
$ws.Cells.Item(406, 2).Value = 'Structure'
$ws.Cells.Item(406, 4).Value = 'Struktur wichtiger als Konvention'
$ws.Cells.Item(406, 7).Value = 'X'
$ws.Cells.Item(406, 13).Value = 'NONAUTOMATIC'

# Row 407
$ws.Cells.Item(407, 2).Value = 'Conditionals'
$ws.Cells.Item(407, 4).Value = 'Nicht triviale Bedingungen in Methoden einkapseln zur besseren Lesbarkeit der Anweisungen'
$ws.Cells.Item(407, 7).Value = 'X'

# Row 408
$ws.Cells.Item(408, 2).Value = 'Conditionals'
$ws.Cells.Item(408, 4).Value = 'Negative Bedingungen vermeiden'
$ws.Cells.Item(408, 7).Value = 'X'

# Row 409
$ws.Cells.Item(409, 2).Value = 'Methods'
$ws.Cells.Item(409, 4).Value = 'Eine Aufgabe pro Methode'
$ws.Cells.Item(409, 7).Value = 'X'
$ws.Cells.Item(409, 13).Value = 'NONAUTOMATIC'

# Row 410
$ws.Cells.Item(410, 2).Value = 'Coupling'
$ws.Cells.Item(410, 4).Value = 'Keine verborgene zeitliche Abhängigkeit (Reihenfolge von Methodenaufrufen)'
$ws.Cells.Item(410, 7).Value = 'X'
$ws.Cells.Item(410, 13).Value = 'NONAUTOMATIC'

# Row 411
$ws.Cells.Item(411, 2).Value = 'Structure'
$ws.Cells.Item(411, 4).Value = 'Keine willkürliche Strukturierung'
$ws.Cells.Item(411, 7).Value = 'X'
$ws.Cells.Item(411, 13).Value = 'NONAUTOMATIC'

# Row 412
$ws.Cells.Item(412, 2).Value = 'Variables'
$ws.Cells.Item(412, 4).Value = 'Grenzbedingungen (z.B. +/- 1) in extra Variable einkapseln'
$ws.Cells.Item(412, 7).Value = 'X'

# Row 413
$ws.Cells.Item(413, 2).Value = 'Structure'
$ws.Cells.Item(413, 3).Value = 'Variables'
$ws.Cells.Item(413, 4).Value = 'Konfigurierbare Konstanten / Variablen auf hoher Abstraktionsebene ansiedeln'
$ws.Cells.Item(413, 7).Value = 'X'
$ws.Cells.Item(413, 13).Value = 'NONAUTOMATIC'

# Row 414
$ws.Cells.Item(414, 2).Value = 'Structure'
$ws.Cells.Item(414, 3).Value = 'Methods'
$ws.Cells.Item(414, 4).Value = 'Transitivie Methodenaufrufe vermeiden - "a.getB().getC()"'
$ws.Cells.Item(414, 7).Value = 'X'

# Row 415
$ws.Cells.Item(415, 2).Value = 'Import'
$ws.Cells.Item(415, 4).Value = 'Wildcard imports nutzen wenn mehr als zwei Klassen aus Package'
$ws.Cells.Item(415, 7).Value = 'X'
$ws.Cells.Item(415, 13).Value = 'Gegenteil von SQ-Regel'

# Row 416
$ws.Cells.Item(416, 2).Value = 'Inheritance'
$ws.Cells.Item(416, 4).Value = 'Keine Konstanten vererben'
$ws.Cells.Item(416, 7).Value = 'X'

# Row 417
$ws.Cells.Item(417, 2).Value = 'Enum'
$ws.Cells.Item(417, 4).Value = 'Enums gegenüber Konstanten bevorzugen'
$ws.Cells.Item(417, 7).Value = 'X'
$ws.Cells.Item(417, 13).Value = 'NONAUTOMATIC'

# --- Step 4: column width adjustments for L (SonarQube-Tags) and M (Clean-Code-Comment) ---
$ws.Columns("L").ColumnWidth = 57.333333333333336
$ws.Columns("M").ColumnWidth = 44

# --- Step 5: update view/selection to match the edited workbook state ---
$ws.Range("D389").Select()
